{"js": "// Replace each two-digit-division answer string with its updated value.\n// All \"old\" strings are unique within the document, so a direct\n// search-and-replace (matchCase, no wildcards) for each pair is safe.\nconst replacements = [\n  [\"39\u00f72=19, 1\", \"67\u00f74=16, 3\"],\n  [\"30\u00f73=10, 0\", \"77\u00f77=11, 0\"],\n  [\"85\u00f73=28, 1\", \"97\u00f74=24, 1\"],\n  [\"93\u00f79=10, 3\", \"49\u00f79=5, 4\"],\n  [\"67\u00f78=8, 3\", \"89\u00f73=29, 2\"],\n  [\"82\u00f73=27, 1\", \"15\u00f74=3, 3\"],\n  [\"65\u00f75=13, 0\", \"82\u00f76=13, 4\"],\n  [\"83\u00f77=11, 6\", \"41\u00f72=20, 1\"],\n  [\"26\u00f79=2, 8\", \"85\u00f77=12, 1\"],\n  [\"57\u00f73=19, 0\", \"10\u00f78=1, 2\"],\n  [\"73\u00f74=18, 1\", \"18\u00f75=3, 3\"],\n  [\"46\u00f74=11, 2\", \"99\u00f75=19, 4\"],\n  [\"12\u00f75=2, 2\", \"66\u00f73=22, 0\"],\n  [\"95\u00f72=47, 1\", \"70\u00f72=35, 0\"],\n  [\"61\u00f73=20, 1\", \"74\u00f72=37, 0\"],\n  [\"17\u00f79=1, 8\", \"95\u00f78=11, 7\"],\n  [\"86\u00f79=9, 5\", \"72\u00f78=9, 0\"],\n  [\"24\u00f77=3, 3\", \"25\u00f75=5, 0\"],\n  [\"75\u00f74=18, 3\", \"39\u00f78=4, 7\"],\n  [\"56\u00f75=11, 1\", \"49\u00f77=7, 0\"],\n  [\"23\u00f79=2, 5\", \"34\u00f79=3, 7\"],\n  [\"93\u00f78=11, 5\", \"58\u00f75=11, 3\"],\n  [\"10\u00f73=3, 1\", \"94\u00f79=10, 4\"],\n  [\"27\u00f72=13, 1\", \"26\u00f72=13, 0\"],\n  [\"49\u00f76=8, 1\", \"78\u00f78=9, 6\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-division answer string with its updated value.\n# All \"old\" strings are unique within the document, so a direct\n# Find/Replace for each pair is safe and won't clobber other cells.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"39\u00f72=19, 1\", \"67\u00f74=16, 3\"),\n    @(\"30\u00f73=10, 0\", \"77\u00f77=11, 0\"),\n    @(\"85\u00f73=28, 1\", \"97\u00f74=24, 1\"),\n    @(\"93\u00f79=10, 3\", \"49\u00f79=5, 4\"),\n    @(\"67\u00f78=8, 3\", \"89\u00f73=29, 2\"),\n    @(\"82\u00f73=27, 1\", \"15\u00f74=3, 3\"),\n    @(\"65\u00f75=13, 0\", \"82\u00f76=13, 4\"),\n    @(\"83\u00f77=11, 6\", \"41\u00f72=20, 1\"),\n    @(\"26\u00f79=2, 8\", \"85\u00f77=12, 1\"),\n    @(\"57\u00f73=19, 0\", \"10\u00f78=1, 2\"),\n    @(\"73\u00f74=18, 1\", \"18\u00f75=3, 3\"),\n    @(\"46\u00f74=11, 2\", \"99\u00f75=19, 4\"),\n    @(\"12\u00f75=2, 2\", \"66\u00f73=22, 0\"),\n    @(\"95\u00f72=47, 1\", \"70\u00f72=35, 0\"),\n    @(\"61\u00f73=20, 1\", \"74\u00f72=37, 0\"),\n    @(\"17\u00f79=1, 8\", \"95\u00f78=11, 7\"),\n    @(\"86\u00f79=9, 5\", \"72\u00f78=9, 0\"),\n    @(\"24\u00f77=3, 3\", \"25\u00f75=5, 0\"),\n    @(\"75\u00f74=18, 3\", \"39\u00f78=4, 7\"),\n    @(\"56\u00f75=11, 1\", \"49\u00f77=7, 0\"),\n    @(\"23\u00f79=2, 5\", \"34\u00f79=3, 7\"),\n    @(\"93\u00f78=11, 5\", \"58\u00f75=11, 3\"),\n    @(\"10\u00f73=3, 1\", \"94\u00f79=10, 4\"),\n    @(\"27\u00f72=13, 1\", \"26\u00f72=13, 0\"),\n    @(\"49\u00f76=8, 1\", \"78\u00f78=9, 6\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n$d.Save()\n"}
